$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'229.54"
$ws.Range("G2").Value = "'4"
$ws.Range("D3").Value = "'22.41"
$ws.Range("G3").Value = "'4"
$ws.Range("D4").Value = "'5.255"
$ws.Range("G4").Value = "'4"
$ws.Range("D5").Value = "'0.05535"
$ws.Range("G5").Value = "'4"
$ws.Range("D6").Value = "'3.381"
$ws.Range("G6").Value = "'4"
$ws.Range("D7").Value = "'6.473"
$ws.Range("G7").Value = "'4"
$ws.Range("D8").Value = "'1.066"
$ws.Range("G8").Value = "'4"
$ws.Range("D9").Value = "'0.7755"
$ws.Range("G9").Value = "'4"
$ws.Range("D10").Value = "'0.1360"
$ws.Range("G10").Value = "'4"
$ws.Range("D11").Value = "'0.07337"
$ws.Range("G11").Value = "'4"
$ws.Range("D12").Value = "'0.03153"
$ws.Range("G12").Value = "'4"
$ws.Range("D13").Value = "'0.02944"
$ws.Range("G13").Value = "'4"
$ws.Range("D14").Value = "'0.09264"
$ws.Range("G14").Value = "'4"
$ws.Range("D15").Value = "'0.001661"
$ws.Range("G15").Value = "'4"
$ws.Range("D16").Value = "'3.249"
$ws.Range("G16").Value = "'4"
$ws.Range("D17").Value = "'0.04774"
$ws.Range("G17").Value = "'4"
$ws.Range("D18").Value = "'0.0005889"
$ws.Range("E18").Value = "17OneONE"
$ws.Range("G18").Value = "'4"
$ws.Range("D19").Value = "'0.006225"
$ws.Range("G19").Value = "'4"
$ws.Range("D20").Value = "'0.005232"
$ws.Range("G20").Value = "'4"
$ws.Range("B21").Value = "BitKan"
$ws.Range("C21").Value = "https://coinranking.com/coin/RDOsLDgvY-AXe+bitkan-kan"
$ws.Range("D21").Value = "'0.001064"
$ws.Range("E21").Value = "20BitKanKAN"
$ws.Range("G21").Value = "'4"
$ws.Range("B22").Value = "NitroEx"
$ws.Range("C22").Value = "https://coinranking.com/coin/8oiZw6gwYhC+nitroex-ntx"
$ws.Range("D22").Value = "'0.0001500"
$ws.Range("E22").Value = "21NitroExNTX"
$ws.Range("G22").Value = "'4"
$ws.Range("B23").Value = "LEO"
$ws.Range("C23").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D23").Value = "'3.940"
$ws.Range("E23").Value = "22LEOLEO"
$ws.Range("G23").Value = "'4"
$ws.Range("B24").Value = "BTSEToken"
$ws.Range("C24").Value = "https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"
$ws.Range("D24").Value = "'2.196"
$ws.Range("E24").Value = "23BTSETokenBTSE"
$ws.Range("G24").Value = "'4"
$ws.Range("B25").Value = "BitpandaEcosystemToken"
$ws.Range("C25").Value = "https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best"
$ws.Range("D25").Value = "'0.3323"
$ws.Range("E25").Value = "24BitpandaEcosystemTokenBEST"
$ws.Range("G25").Value = "'4"
$ws.Range("B26").Value = "ProBitToken"
$ws.Range("C26").Value = "https://coinranking.com/coin/lQP4d6T2+probittoken-prob"
$ws.Range("D26").Value = "'0.1243"
$ws.Range("E26").Value = "25ProBitTokenPROB"
$ws.Range("G26").Value = "'4"
$ws.Range("B27").Value = "UpBots"
$ws.Range("C27").Value = "https://coinranking.com/coin/m5ozaAIK6+upbots-ubxt"
$ws.Range("D27").Value = "'0.0004999"
$ws.Range("E27").Value = "26UpBotsUBXTBestin24h"
$ws.Range("G27").Value = "'4"
$ws.Range("G28").Value = "'4"
$ws.Range("G29").Value = "'4"
$ws.Range("G30").Value = "'4"
$ws.Range("G31").Value = "'4"
$ws.Range("G32").Value = "'4"
$ws.Range("G33").Value = "'4"
$ws.Range("G34").Value = "'4"
$ws.Range("G35").Value = "'4"
$ws.Range("G36").Value = "'4"
$ws.Range("G37").Value = "'4"
$ws.Range("G38").Value = "'4"
$ws.Range("G39").Value = "'4"
$ws.Range("D40").Value = "'0.03955"
$ws.Range("G40").Value = "'4"
$ws.Range("D41").Value = "'0.007130"
$ws.Range("G41").Value = "'4"
$ws.Range("B42").Value = "BKEXToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
$ws.Range("D42").Value = "'0.1037"
$ws.Range("E42").Value = "41BKEXTokenBKK"
$ws.Range("G42").Value = "'4"
$ws.Range("B43").Value = "CEJI"
$ws.Range("C43").Value = "https://coinranking.com/coin/SbKjCVJCh+ceji-ceji"
$ws.Range("D43").Value = "'0.002679"
$ws.Range("E43").Value = "42CEJICEJI"
$ws.Range("G43").Value = "'4"
$ws.Range("D44").Value = "'0.008614"
$ws.Range("G44").Value = "'4"
$ws.Range("D45").Value = "'0.00005443"
$ws.Range("G45").Value = "'4"
$ws.Range("G46").Value = "'4"
$ws.Range("D47").Value = "'0.7850"
$ws.Range("G47").Value = "'4"
$ws.Range("D48").Value = "'0.04011"
$ws.Range("E48").Value = "47BOLOBOLOWorstin24h"
$ws.Range("G48").Value = "'4"
$ws.Range("G49").Value = "'4"
$ws.Range("G50").Value = "'4"
$ws.Range("G51").Value = "'4"
